$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the value in C5 (G5's shared formula recalculates from this)
$ws.Range("C5").Value = 37.9

# Extend the AVERAGE formulas in C21/D21 to also include row 5 (G5/H5)
$ws.Range("C21").Formula = "=AVERAGE(G4,G6,G7,G8,G5)"
$ws.Range("D21").Formula = "=AVERAGE(H4,H6,H7,H8,H5)"

# Match the final view/selection state (D21 active)
$ws.Activate()
$ws.Range("D21").Select()
